# Add a new test row to the "Location" sheet for the "Home.Group.All"
# click-group link (commit: "Add test on click group All").
#
# The Location sheet lists Key / Type / Value rows (row 1 = headers).
# Existing rows 3-6 describe the "Home.Group.<Era>" LinkText rows; this
# adds the analogous row for the "All" group link at the end of the
# table (row 8), reusing the same alternating row style (copied from
# row 6, which already carries the even-row style used by rows 2/4/6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Location")

# Copy formatting (styles) of the last existing data row down into the
# new row so the new row matches the established look (fill/border/etc.)
$ws.Range("A6:G6").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)

# Populate the new row's content.
$ws.Range("A8").Value = "Home.Group.All"
$ws.Range("B8").Value = "LinkText"
$ws.Range("C8").Value = "All"

# Match the row height used by the other data rows in this table.
$ws.Rows.Item(8).RowHeight = 20.25
